$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update temperature (F) and windspeed (G) columns with corrected export data
$ws.Range("F2").Value = 21.4
$ws.Range("G2").Value = 13.8

$ws.Range("F3").Value = 29.3
$ws.Range("G3").Value = 14.8

$ws.Range("F4").Value = 30.1
$ws.Range("G4").Value = 5.9

$ws.Range("F5").Value = 12.8
$ws.Range("G5").Value = 8.4

$ws.Range("F6").Value = 24
$ws.Range("G6").Value = 12.2

$ws.Range("F7").Value = 28.9
$ws.Range("G7").Value = 7.1

$ws.Range("F8").Value = 17.2
$ws.Range("G8").Value = 6.8

$ws.Range("F9").Value = 32.9
$ws.Range("G9").Value = 7.9

$ws.Range("F10").Value = 36.3
$ws.Range("G10").Value = 21.1

$ws.Range("F11").Value = 30.3
$ws.Range("G11").Value = 31.3

$ws.Range("F12").Value = 10.9
$ws.Range("G12").Value = 16

$ws.Range("F13").Value = -0
$ws.Range("G13").Value = 7.2

$ws.Range("F14").Value = 7.5
$ws.Range("G14").Value = 6.9

$ws.Range("F15").Value = 12.6
$ws.Range("G15").Value = 37.5

$ws.Range("F16").Value = 13.2
$ws.Range("G16").Value = 23.8

$ws.Range("F17").Value = 26.4
$ws.Range("G17").Value = 33.7

$ws.Range("F18").Value = 7.8
$ws.Range("G18").Value = 12.2

$ws.Range("F19").Value = 3.3
$ws.Range("G19").Value = 1.3

$ws.Range("F20").Value = 1.7
$ws.Range("G20").Value = 0.7

$ws.Range("F21").Value = 9.699999999999999
$ws.Range("G21").Value = 2.1

# Remove the currency_group column entirely (money graph data no longer needed here)
$ws.Columns("H").Delete()
